$wb = $excel.ActiveWorkbook

# Rename sheets (formulas/chart refs that point at these sheets follow the rename automatically)
$wb.Worksheets.Item("l11_gbalo").Name = "t0"
$wb.Worksheets.Item("Sheet1").Name = "t1"
$wb.Worksheets.Item("Sheet2").Name = "t2"
$wb.Worksheets.Item("Sheet3").Name = "t3"

# Update the selection on the active sheet (t3, formerly Sheet3) from L17 to E15
$ws3 = $wb.Worksheets.Item("t3")
$ws3.Activate()
$ws3.Range("E15").Select()
